$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3007.762
$ws.Range("I62").Value = 2103.6875
$ws.Range("J62").Value = 5900.8
$ws.Range("K62").Value = 2103.6875
$ws.Range("L62").Value = 5900.8
$ws.Range("M62").Value = -1479.6875
$ws.Range("N62").Value = -7148.8
$ws.Range("H65").Value = 3007.762
$ws.Range("I65").Value = 2103.6875
$ws.Range("J65").Value = 5900.8
$ws.Range("K65").Value = 10518.4375
$ws.Range("L65").Value = 29504
$ws.Range("M65").Value = -7398.4375
$ws.Range("N65").Value = -35744
$ws.Range("H113").Value = 2767.4
$ws.Range("J113").Value = 2771.75
$ws.Range("L113").Value = 2771.75
$ws.Range("N113").Value = -9279.75
$ws.Range("H132").Value = 1919.625
$ws.Range("I132").Value = 1919.625
$ws.Range("K132").Value = 5758.875
$ws.Range("M132").Value = -3228.875
$ws.Range("H135").Value = 150001820
$ws.Range("I135").Value = 100000340
$ws.Range("J135").Value = 200003300
$ws.Range("K135").Value = 900003060
$ws.Range("L135").Value = 1800029700
$ws.Range("M135").Value = -900000525
$ws.Range("N135").Value = -1800034770
$ws.Range("H137").Value = 619796
$ws.Range("I137").Value = 2213.037
$ws.Range("J137").Value = 1237379
$ws.Range("K137").Value = 6639.110999999999
$ws.Range("L137").Value = 3712137
$ws.Range("M137").Value = -4089.110999999999
$ws.Range("N137").Value = -3717237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15923.014
$ws.Range("I32").Value = 18326.87
$ws.Range("J32").Value = 3703.4167
$ws.Range("K32").Value = 18326.87
$ws.Range("L32").Value = 3703.4167
$ws.Range("M32").Value = -18039.87
$ws.Range("N32").Value = -4277.4167
$ws.Range("H129").Value = 39644.5
$ws.Range("J129").Value = 39644.5
$ws.Range("L129").Value = 39644.5
$ws.Range("N129").Value = -49644.5
$ws.Range("H132").Value = 1927.6
$ws.Range("I132").Value = 1539
$ws.Range("K132").Value = 4617
$ws.Range("M132").Value = -2087

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 80000
$ws.Range("J59").Value = 80000
$ws.Range("L59").Value = 80000
$ws.Range("N59").Value = -81694
$ws.Range("H99").Value = 2548.25
$ws.Range("I99").Value = 1269.7273
$ws.Range("J99").Value = 4110.8887
$ws.Range("K99").Value = 1269.7273
$ws.Range("L99").Value = 4110.8887
$ws.Range("M99").Value = 228.2727
$ws.Range("N99").Value = -7106.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3364.4883
$ws.Range("I132").Value = 2447.4285
$ws.Range("J132").Value = 5076.3335
$ws.Range("K132").Value = 7342.2855
$ws.Range("L132").Value = 15229.0005
$ws.Range("M132").Value = -4812.2855
$ws.Range("N132").Value = -20289.0005
$ws.Range("H133").Value = 23774
$ws.Range("J133").Value = 37400
$ws.Range("L133").Value = 37400
$ws.Range("N133").Value = -42460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6670978.5
$ws.Range("I5").Value = 338.34784
$ws.Range("J5").Value = 83383336
$ws.Range("K5").Value = 1015.04352
$ws.Range("L5").Value = 250150008
$ws.Range("M5").Value = -903.0435200000001
$ws.Range("N5").Value = -250150232
$ws.Range("H62").Value = 3400
$ws.Range("J62").Value = 3400
$ws.Range("L62").Value = 10200
$ws.Range("N62").Value = -11572
$ws.Range("H65").Value = 3400
$ws.Range("J65").Value = 3400
$ws.Range("L65").Value = 30600
$ws.Range("N65").Value = -37464
$ws.Range("H68").Value = 2593.2307
$ws.Range("I68").Value = 1106.3
$ws.Range("J68").Value = 5248.4644
$ws.Range("K68").Value = 3318.9
$ws.Range("L68").Value = 15745.3932
$ws.Range("M68").Value = -2507.9
$ws.Range("N68").Value = -17367.3932
$ws.Range("H71").Value = 2593.2307
$ws.Range("I71").Value = 1106.3
$ws.Range("J71").Value = 5248.4644
$ws.Range("K71").Value = 9956.699999999999
$ws.Range("L71").Value = 47236.1796
$ws.Range("M71").Value = -5900.699999999999
$ws.Range("N71").Value = -55348.1796
$ws.Range("H107").Value = 462.875
$ws.Range("I107").Value = 273.11667
$ws.Range("K107").Value = 819.35001
$ws.Range("M107").Value = 1100.64999
$ws.Range("H109").Value = 884.55554
$ws.Range("J109").Value = 3030
$ws.Range("L109").Value = 9090
$ws.Range("N109").Value = -11170
$ws.Range("H112").Value = 1613.5
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H120").Value = 10400
$ws.Range("I120").Value = 10500
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 31500
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -26662
$ws.Range("N120").Value = -39676
$ws.Range("H121").Value = 533
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 533
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 1599
$ws.Range("N121").Value = -4219
$ws.Range("H122").Value = 877.2857
$ws.Range("I122").Value = 362.875
$ws.Range("K122").Value = 3265.875
$ws.Range("M122").Value = -815.875
$ws.Range("H125").Value = 2892
$ws.Range("J125").Value = 6500
$ws.Range("L125").Value = 19500
$ws.Range("N125").Value = -29340
$ws.Range("H131").Value = 1117.5135
$ws.Range("I131").Value = 1332.6
$ws.Range("J131").Value = 970.86365
$ws.Range("K131").Value = 3997.8
$ws.Range("L131").Value = 2912.59095
$ws.Range("M131").Value = 1042.2
$ws.Range("N131").Value = -12992.59095
$ws.Range("H134").Value = 4381.483
$ws.Range("I134").Value = 4775.154
$ws.Range("J134").Value = 4061.625
$ws.Range("K134").Value = 14325.462
$ws.Range("L134").Value = 12184.875
$ws.Range("M134").Value = -9255.462000000001
$ws.Range("N134").Value = -22324.875
$ws.Range("H135").Value = 6670978.5
$ws.Range("I135").Value = 338.34784
$ws.Range("J135").Value = 83383336
$ws.Range("K135").Value = 3045.13056
$ws.Range("L135").Value = 750450024
$ws.Range("M135").Value = -510.1305600000001
$ws.Range("N135").Value = -750455094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5771.5713
$ws.Range("I107").Value = 12242.167
$ws.Range("J107").Value = 918.625
$ws.Range("K107").Value = 12242.167
$ws.Range("L107").Value = 918.625
$ws.Range("M107").Value = -10322.167
$ws.Range("N107").Value = -4758.625
$ws.Range("H113").Value = 2085.5652
$ws.Range("I113").Value = 1937.2
$ws.Range("J113").Value = 2363.75
$ws.Range("K113").Value = 1937.2
$ws.Range("L113").Value = 2363.75
$ws.Range("M113").Value = 232.8
$ws.Range("N113").Value = -6703.75
$ws.Range("H126").Value = 2508
$ws.Range("I126").Value = 1814.9333
$ws.Range("J126").Value = 3307.6924
$ws.Range("K126").Value = 5444.7999
$ws.Range("L126").Value = 9923.0772
$ws.Range("M126").Value = -2974.7999
$ws.Range("N126").Value = -14863.0772
$ws.Range("H133").Value = 63589.332
$ws.Range("J133").Value = 63589.332
$ws.Range("L133").Value = 63589.332
$ws.Range("N133").Value = -73709.33199999999
$ws.Range("H136").Value = 26799.4
$ws.Range("J136").Value = 26799.4
$ws.Range("L136").Value = 80398.20000000001
$ws.Range("N136").Value = -85498.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2947.7646
$ws.Range("I7").Value = 3246.7
$ws.Range("J7").Value = 2520.7144
$ws.Range("K7").Value = 3246.7
$ws.Range("L7").Value = 2520.7144
$ws.Range("M7").Value = -3134.7
$ws.Range("N7").Value = -2744.7144
$ws.Range("H93").Value = 910.8182
$ws.Range("I93").Value = 779.44446
$ws.Range("J93").Value = 1502
$ws.Range("K93").Value = 779.44446
$ws.Range("L93").Value = 1502
$ws.Range("M93").Value = 468.55554
$ws.Range("N93").Value = -3998
$ws.Range("H126").Value = 2947.7646
$ws.Range("I126").Value = 3246.7
$ws.Range("J126").Value = 2520.7144
$ws.Range("K126").Value = 9740.099999999999
$ws.Range("L126").Value = 7562.1432
$ws.Range("M126").Value = -7270.099999999999
$ws.Range("N126").Value = -12502.1432
$ws.Range("H133").Value = 46246.75
$ws.Range("J133").Value = 46246.75
$ws.Range("L133").Value = 46246.75
$ws.Range("N133").Value = -51306.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1792.6111
$ws.Range("I126").Value = 1792.6111
$ws.Range("K126").Value = 5377.8333
$ws.Range("M126").Value = -2907.8333
